{"js": "// The diff re-saves the run \"Obeo's website\" (inside the w:hyperlink) with\n// its run properties re-serialized (POI 4.1.0 -> 5.2.3 upgrade): the\n// boolean toggles b/i/strike go from explicit true/false to the equivalent\n// on/off form, i.e. the *logical* formatting is unchanged - bold stays on,\n// italic stays off, strikethrough stays off. Re-apply the same formatting\n// explicitly on that run so it is (re)written.\n\nconst results = context.document.body.search(\"Obeo's website\", { matchCase: true, matchWholeWord: false });\nresults.load(\"font\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find run with text \\\"Obeo's website\\\"\");\n}\n\nconst font = results.items[0].font;\n\n// Re-assert the same logical values (true/false <-> on/off are equivalent\n// in OOXML); this is what actually changed serialization-wise in the diff.\nfont.bold = true;\nfont.italic = false;\nfont.strikeThrough = false;\n\nawait context.sync();\n", "ps1": "# The diff re-saves the run \"Obeo's website\" (inside the w:hyperlink) with\n# its run properties re-serialized (POI 4.1.0 -> 5.2.3 upgrade): the\n# boolean toggles b/i/strike go from explicit true/false to the equivalent\n# on/off form, i.e. the *logical* formatting is unchanged - bold stays on,\n# italic stays off, strikethrough stays off. Re-apply the same formatting\n# explicitly on that run so it is (re)written.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Obeo's website\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found -and $find.Found) {\n    $rng.Font.Bold = $true\n    $rng.Font.Italic = $false\n    $rng.Font.StrikeThrough = $false\n}\n"}
